$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.773.30"
$ws.Range("D2").Style = $s
$ws.Range("E2").Value = "  -4.37%  "

$s = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.497.79"
$ws.Range("D3").Style = $s
$ws.Range("E3").Value = "  -5.29%  "

$ws.Range("E4").Value = "  -0.02%  "

$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.71"
$ws.Range("D5").Style = $s
$ws.Range("E5").Value = "  -1.80%  "

$s = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.92"
$ws.Range("D6").Style = $s
$ws.Range("E6").Value = "  -3.52%  "

$ws.Range("E7").Value = "  -0.47%  "

$s = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.490.01"
$ws.Range("D8").Style = $s
$ws.Range("E8").Value = "  -5.27%  "

$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("E10").Value = "  -7.48%  "

$s = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.64"
$ws.Range("D11").Style = $s
$ws.Range("E11").Value = "  +9.23%  "

$ws.Range("E12").Value = "  -2.40%  "

$s = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.11"
$ws.Range("D13").Style = $s
$ws.Range("E13").Value = "  -5.95%  "

$ws.Range("E14").Value = "  -3.91%  "

$s = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "671.76"
$ws.Range("D15").Style = $s
$ws.Range("E15").Value = "  -1.96%  "

$s = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.059.20"
$ws.Range("D16").Style = $s
$ws.Range("E16").Value = "  -5.35%  "

$s = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.81"
$ws.Range("D17").Style = $s
$ws.Range("E17").Value = "  -1.83%  "

$s = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.500.14"
$ws.Range("D18").Style = $s
$ws.Range("E18").Value = "  -5.17%  "

$s = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.733.78"
$ws.Range("D19").Style = $s
$ws.Range("E19").Value = "  -4.59%  "

$ws.Range("E20").Value = "  -1.72%  "

$ws.Range("E21").Value = "  -4.48%  "

$s = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.17"
$ws.Range("D22").Style = $s
$ws.Range("E22").Value = "  -4.35%  "

$ws.Range("E23").Value = "  -4.32%  "

$s = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.24"
$ws.Range("D24").Style = $s
$ws.Range("E24").Value = "  -8.80%  "

$s = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.06"
$ws.Range("D25").Style = $s
$ws.Range("E25").Value = "  -5.49%  "

$s = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.86"
$ws.Range("D26").Style = $s
$ws.Range("E26").Value = "  -4.57%  "

$s = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.80"
$ws.Range("D27").Style = $s
$ws.Range("E27").Value = "  -1.09%  "

$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("E29").Value = "  -7.66%  "

$s = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.40"
$ws.Range("D30").Style = $s
$ws.Range("E30").Value = "  -8.16%  "

$s = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.84"
$ws.Range("D31").Style = $s
$ws.Range("E31").Value = "  -7.74%  "

$s = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.69"
$ws.Range("D32").Style = $s
$ws.Range("E32").Value = "  -5.61%  "

$s = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.19"
$ws.Range("D33").Style = $s
$ws.Range("E33").Value = "  -8.60%  "

$ws.Range("E34").Value = "  -1.74%  "

$ws.Range("E35").Value = "  -5.92%  "

$s = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "577.97"
$ws.Range("D36").Style = $s
$ws.Range("E36").Value = "  -0.20%  "

$s = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.91"
$ws.Range("D37").Style = $s
$ws.Range("E37").Value = "  -3.77%  "

$s = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.57"
$ws.Range("D38").Style = $s
$ws.Range("E38").Value = "  -15.16%  "

$ws.Range("E39").Value = "  -4.55%  "

$s = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.84"
$ws.Range("D40").Style = $s
$ws.Range("E40").Value = "  -5.08%  "

$s = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = $s
$ws.Range("E41").Value = "  +0.09%  "

$s = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.136"
$ws.Range("D42").Style = $s
$ws.Range("E42").Value = "  -5.80%  "

$ws.Range("E43").Value = "  -3.77%  "

$s = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0436"
$ws.Range("D44").Style = $s
$ws.Range("E44").Value = "  -5.72%  "

$s = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.410.86"
$ws.Range("D45").Style = $s
$ws.Range("E45").Value = "  -9.38%  "

$s = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "33.27"
$ws.Range("D46").Style = $s
$ws.Range("E46").Value = "  -6.58%  "

$ws.Range("E47").Value = "  -9.55%  "

$ws.Range("E48").Value = "  +0.92%  "

$ws.Range("E49").Value = "  -7.68%  "

$ws.Range("E50").Value = "  -0.91%  "

$s = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.98"
$ws.Range("D51").Style = $s
$ws.Range("E51").Value = "  -1.48%  "
